# Actualización automatica mar abr  6 17:33:39 CEST 2021
#
# Fixes the metadata annotation block (rows 3-5) of the report sheet:
#  - column A (the "orden" column) gets its missing sdmx/measure metadata
#  - column D ("siglas") gets corrected metadata (was tagged as a
#    dimension/skos concept, should be a measure/xsd:string)
#  - the stray trailing note in D6 is removed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "iaest-measure:orden"
$ws.Range("A4").Value = "medida"
$ws.Range("A5").Value = "xsd:int"

$ws.Range("D3").Value = "iaest-measure:siglas"
$ws.Range("D4").Value = "medida"
$ws.Range("D5").Value = "xsd:string"

$ws.Rows(6).Delete()
